$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.480.87'
$ws.Range('E2').Value = '  +9.52%  '
$ws.Range('D3').Value = '3.330.22'
$ws.Range('E3').Value = '  +6.34%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '217.27'
$ws.Range('E5').Value = '  +5.97%  '
$ws.Range('D6').Value = '650.09'
$ws.Range('E6').Value = '  +4.51%  '
$ws.Range('D7').Value = '0.396'
$ws.Range('E7').Value = '  +41.29%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  +5.08%  '
$ws.Range('D10').Value = '3.326.36'
$ws.Range('E10').Value = '  +6.25%  '
$ws.Range('D11').Value = '0.585'
$ws.Range('E11').Value = '  +2.23%  '
$ws.Range('D12').Value = '0.0000289'
$ws.Range('E12').Value = '  +16.07%  '
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '35.33'
$ws.Range('E14').Value = '  +13.56%  '
$ws.Range('D15').Value = '3.944.14'
$ws.Range('E15').Value = '  +6.09%  '
$ws.Range('D16').Value = '5.50'
$ws.Range('E16').Value = '  +4.99%  '
$ws.Range('D17').Value = '88.419.80'
$ws.Range('E17').Value = '  +9.42%  '
$ws.Range('D18').Value = '3.334.00'
$ws.Range('E18').Value = '  +5.95%  '
$ws.Range('D19').Value = '14.60'
$ws.Range('E19').Value = '  +5.38%  '
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').Value = '3.11'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '9.67'
$ws.Range('E21').Value = '  +8.68%  '
$ws.Range('D22').Value = '454.72'
$ws.Range('E22').Value = '  +5.84%  '
$ws.Range('D23').Value = '5.45'
$ws.Range('E23').Value = '  +7.91%  '
$ws.Range('E24').Value = '  +4.55%  '
$ws.Range('D25').Value = '5.52'
$ws.Range('E25').Value = '  +7.84%  '
$ws.Range('D26').Value = '12.80'
$ws.Range('E26').Value = '  +18.43%  '
$ws.Range('D27').Value = '3.514.36'
$ws.Range('E27').Value = '  +6.10%  '
$ws.Range('D28').Value = '78.52'
$ws.Range('E28').Value = '  +3.92%  '
$ws.Range('D29').Value = '0.212'
$ws.Range('E29').Value = '  +42.95%  '
$ws.Range('E30').Value = '  +10.84%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  +5.05%  '
$ws.Range('D33').Value = '590.46'
$ws.Range('E33').Value = '  +6.99%  '
$ws.Range('D34').Value = '1.59'
$ws.Range('E34').Value = '  +8.57%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  +6.48%  '
$ws.Range('D37').Value = '7.13'
$ws.Range('E37').Value = '  +21.51%  '
$ws.Range('D38').Value = '0.142'
$ws.Range('E38').Value = '  -4.84%  '
$ws.Range('D39').Value = '23.08'
$ws.Range('E39').Value = '  +2.45%  '
$ws.Range('D40').Value = '2.14'
$ws.Range('E40').Value = '  +8.82%  '
$ws.Range('D41').Value = '0.420'
$ws.Range('E41').Value = '  +4.23%  '
$ws.Range('D42').Value = '21.84'
$ws.Range('E42').Value = '  +5.54%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '3.15'
$ws.Range('E44').Value = '  +5.34%  '
$ws.Range('D45').Value = '158.19'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').Value = '  +9.67%  '
$ws.Range('D48').Value = '187.85'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('D49').Value = '46.00'
$ws.Range('E49').Value = '  +5.28%  '
$ws.Range('E50').Value = '  +5.88%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '0.659'
$ws.Range('E51').Value = '  +6.02%  '
